$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "Date.x"
$ws.Range("D37").Value = "Date"

$ws.Range("A38").Value = "Time.x"
$ws.Range("D38").Value = "Time"

$ws.Range("A39").Select()
